# Adds new progress as of 04-Nov-2025:
#   - Column H ("PERIOD TO EXPIRE") decreases by 1 day for every data row.
#   - Column I ("LAST UPDATE") moves from 03-Nov-2025 to 04-Nov-2025.
# Applies to the "Training Dashboard" sheet, data rows 3-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$firstRow = 3
$lastRow = 27

for ($row = $firstRow; $row -le $lastRow; $row++) {

    # --- Column H: PERIOD TO EXPIRE, numeric countdown, decrement by 1 ---
    $hCell = $ws.Cells.Item($row, 8)
    $hVal = $hCell.Value2
    if ($hVal -ne $null) {
        $hCell.Value = $hVal - 1
    }

    # --- Column I: LAST UPDATE, stored as literal text "dd-mmm-yyyy" ---
    $iCell = $ws.Cells.Item($row, 9)
    if ($iCell.Value2 -ne $null) {
        # Force text storage first so Excel doesn't reinterpret the literal
        # "04-Nov-2025" string as a date serial number.
        $iCell.NumberFormat = "@"
        $iCell.Value = "04-Nov-2025"

        # Restore the original cell formatting/style (border, alignment,
        # General number format) by pasting formats from an untouched
        # neighbor cell on the same row that shares the same style.
        $ws.Cells.Item($row, 10).Copy() | Out-Null
        $iCell.PasteSpecial(-4122) | Out-Null
    }
}

$excel.CutCopyMode = 0
